$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.505.37"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.618.08"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.30"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.526"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.80"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.845.40"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "1.614.92"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.48"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "27.520.65"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.45"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.92"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +6.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.95"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.83"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "1.443.42"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.05"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.942"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.560"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.863"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.18"
$ws.Range("E41").Value = "  +6.60%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.39"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "1.756.49"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.22"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0988"
$ws.Range("E51").Value = "  +0.18%  "
